$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day of price data (2022-02-04) is inserted at row 46, pushing the
# existing rows 46:69 down to 47:70.
$ws.Rows("46:46").Insert()

$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44596
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = 100112038
$ws.Range("G46").Value = "Cebollín baby"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 300
$ws.Range("K46").Value = 2800
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = 2900
$ws.Range("N46").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O46").Value = "Región de Arica y Parinacota"
$ws.Range("P46").Value = 1450
$ws.Range("Q46").Value = 2
$ws.Range("R46").Value = "Hortaliza"
